$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "<roil>"
$ws.Range("C2").Value = 37

# Row 3
$ws.Range("C3").Value = 37

# Row 4
$ws.Range("B4").Value = "<her>"
$ws.Range("C4").Value = 36

# Row 5
$ws.Range("C5").Value = 38

# Row 6
$ws.Range("C6").Value = 43

# Row 7
$ws.Range("C7").Value = 37

# Row 8
$ws.Range("B8").Value = "<tho>"
$ws.Range("C8").Value = 34

# Row 9
$ws.Range("C9").Value = 32

# Row 10
$ws.Range("C10").Value = 42

# Row 11
$ws.Range("B11").Value = "<an>"
$ws.Range("C11").Value = 35

# Row 12
$ws.Range("B12").Value = "<eight>"
$ws.Range("C12").Value = 33

# Row 13
$ws.Range("B13").Value = "<in>"
$ws.Range("C13").Value = 44

# Row 14
$ws.Range("C14").Value = 35

# Row 15
$ws.Range("B15").Value = "<of>"
$ws.Range("C15").Value = 40

# Row 16
$ws.Range("B16").Value = "<frim>"
$ws.Range("C16").Value = 20
